# Adds the "NOISE MODE" screen translation rows (355-370) to the
# "Translation" sheet, per the commit:
#   aggiunta la schermata NOISE_MODE con logica, inserito il controllo sul
#   valore massimo degli indicatori e controllo sull'indicatore Fuel Level
#
# Each row holds: Text ID (B), Typography Name (C), Alignment (D),
# GB text (E), Language/LTR marker (F).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

$rows = @(
    @{ Row = 355; B = "SingleUseId379"; C = "Little";       D = "Center"; E = "NOISE MODE"; F = "LTR" },
    @{ Row = 356; B = "SingleUseId380"; C = "LittleMedium";  D = "Center"; E = "RPM:";       F = "LTR" },
    @{ Row = 357; B = "SingleUseId381"; C = "Medium";        D = "Right";  E = "<value>";    F = "LTR" },
    @{ Row = 358; B = "SingleUseId383"; C = "LittleMedium";  D = "Center"; E = "TPS:";       F = "LTR" },
    @{ Row = 359; B = "SingleUseId384"; C = "Medium";        D = "Right";  E = "<value>";    F = "LTR" },
    @{ Row = 360; B = "SingleUseId385"; C = "Medium";        D = "Left";   E = "0000";       F = "LTR" },
    @{ Row = 361; B = "SingleUseId386"; C = "LittleMedium";  D = "Center"; E = "TRMC 2:";    F = "LTR" },
    @{ Row = 362; B = "SingleUseId387"; C = "Medium";        D = "Right";  E = "<value>";    F = "LTR" },
    @{ Row = 363; B = "SingleUseId388"; C = "Medium";        D = "Left";   E = "0000";       F = "LTR" },
    @{ Row = 364; B = "SingleUseId389"; C = "Bold";          D = "Center"; E = "TRMC:";      F = "LTR" },
    @{ Row = 365; B = "SingleUseId390"; C = "LittleMedium";  D = "Center"; E = "TRMC 1:";    F = "LTR" },
    @{ Row = 366; B = "SingleUseId391"; C = "Medium";        D = "Right";  E = "<value>";    F = "LTR" },
    @{ Row = 367; B = "SingleUseId392"; C = "Medium";        D = "Left";   E = "0000";       F = "LTR" },
    @{ Row = 368; B = "SingleUseId394"; C = "Large";         D = "Center"; E = "<value>";    F = "LTR" },
    @{ Row = 369; B = "SingleUseId396"; C = "Medium";        D = "Left";   E = "0000";       F = "LTR" },
    @{ Row = 370; B = "SingleUseId397"; C = "Large";         D = "Left";   E = "N";          F = "LTR" }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("B$n").Value = $r.B
    $ws.Range("C$n").Value = $r.C
    $ws.Range("D$n").Value = $r.D
    $ws.Range("E$n").Value = $r.E
    $ws.Range("F$n").Value = $r.F
}
